$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (cellRef, action, value)
# action "SET" assigns a numeric value; "DELETE" clears the cell so the
# cell node is dropped entirely from the saved OOXML (matches cells that
# disappear from a row in the source diff).
$changes = @{}
$changes["ALC"] = @(
    ,@("H19", "SET", 892.13513)
    ,@("I19", "SET", 782)
    ,@("J19", "SET", 976.0476)
    ,@("K19", "SET", 782)
    ,@("L19", "SET", 976.0476)
    ,@("M19", "SET", -607)
    ,@("N19", "SET", -1326.0476)
    ,@("H74", "SET", 4750)
    ,@("I74", "SET", 0)
    ,@("J74", "SET", 4750)
    ,@("K74", "SET", 0)
    ,@("L74", "DELETE", $null)
    ,@("M74", "SET", 4750)
    ,@("N74", "SET", -6622)
    ,@("H77", "SET", 4750)
    ,@("I77", "SET", 0)
    ,@("J77", "SET", 4750)
    ,@("K77", "SET", 0)
    ,@("L77", "DELETE", $null)
    ,@("M77", "SET", 23750)
    ,@("N77", "SET", -33110)
    ,@("H98", "SET", 1605.3529)
    ,@("I98", "SET", 1072.091)
    ,@("J98", "SET", 2583)
    ,@("K98", "SET", 1072.091)
    ,@("L98", "SET", 2583)
    ,@("M98", "SET", 425.9090000000001)
    ,@("N98", "SET", -5579)
    ,@("H111", "SET", 2646.1924)
    ,@("I111", "SET", 2075.3125)
    ,@("J111", "SET", 3559.6)
    ,@("K111", "SET", 6225.9375)
    ,@("L111", "SET", 10678.8)
    ,@("M111", "SET", -3158.9375)
    ,@("N111", "SET", -16812.8)
    ,@("H112", "SET", 1169.4117)
    ,@("I112", "SET", 586.6667)
    ,@("J112", "SET", 1294.2858)
    ,@("K112", "SET", 1760.0001)
    ,@("L112", "SET", 3882.8574)
    ,@("M112", "SET", -652.0001)
    ,@("N112", "SET", -6098.857400000001)
    ,@("H122", "SET", 1605.3529)
    ,@("I122", "SET", 1072.091)
    ,@("J122", "SET", 2583)
    ,@("K122", "SET", 3216.273)
    ,@("L122", "SET", 7749)
    ,@("M122", "SET", -766.2729999999997)
    ,@("N122", "SET", -12649)
    ,@("H124", "SET", 23636.363)
    ,@("J124", "SET", 23636.363)
    ,@("L124", "SET", 23636.363)
    ,@("N124", "SET", -33456.363)
    ,@("H126", "SET", 23333.334)
    ,@("J126", "SET", 23333.334)
    ,@("L126", "SET", 23333.334)
    ,@("N126", "SET", -33213.334)
    ,@("H129", "SET", 809.8182)
    ,@("I129", "SET", 518)
    ,@("J129", "SET", 1160)
    ,@("K129", "SET", 1554)
    ,@("L129", "SET", 3480)
    ,@("M129", "SET", 3446)
    ,@("N129", "SET", -13480)
    ,@("H130", "SET", 26250)
    ,@("J130", "SET", 26250)
    ,@("L130", "SET", 26250)
    ,@("N130", "SET", -36290)
    ,@("H137", "SET", 1313.125)
    ,@("I137", "SET", 879.5357)
    ,@("J137", "SET", 4348.25)
    ,@("K137", "SET", 2638.6071)
    ,@("L137", "SET", 13044.75)
    ,@("M137", "SET", -88.60710000000017)
    ,@("N137", "SET", -18144.75)
)
$changes["BSM"] = @(
    ,@("H70", "SET", 105900)
    ,@("J70", "SET", 105900)
    ,@("L70", "SET", 105900)
    ,@("N70", "SET", -106486)
    ,@("H73", "SET", 105900)
    ,@("J73", "SET", 105900)
    ,@("L73", "SET", 105900)
    ,@("N73", "SET", -107928)
)
$changes["CRP"] = @(
    ,@("H80", "SET", 35000)
    ,@("J80", "SET", 35000)
    ,@("L80", "SET", 35000)
    ,@("N80", "SET", -37246)
    ,@("H83", "SET", 35000)
    ,@("J83", "SET", 35000)
    ,@("L83", "SET", 105000)
    ,@("N83", "SET", -116232)
    ,@("H129", "SET", 30371.5)
    ,@("J129", "SET", 30371.5)
    ,@("L129", "SET", 30371.5)
    ,@("N129", "SET", -40371.5)
    ,@("H134", "SET", 1641.0857)
    ,@("I134", "SET", 1685.8572)
    ,@("J134", "SET", 1462)
    ,@("K134", "SET", 5057.571599999999)
    ,@("L134", "SET", 4386)
    ,@("M134", "SET", -2522.571599999999)
    ,@("N134", "SET", -9456)
)
$changes["CUL"] = @(
    ,@("H22", "SET", 4100)
    ,@("J22", "SET", 4100)
    ,@("L22", "SET", 12300)
    ,@("N22", "SET", -12638)
    ,@("H27", "SET", 4100)
    ,@("J27", "SET", 4100)
    ,@("L27", "SET", 12300)
    ,@("N27", "SET", -12504)
    ,@("H34", "SET", 481.83334)
    ,@("I34", "SET", 360.6)
    ,@("J34", "SET", 568.4286)
    ,@("K34", "SET", 1081.8)
    ,@("L34", "SET", 1705.2858)
    ,@("M34", "SET", -997.8000000000002)
    ,@("N34", "SET", -1873.2858)
    ,@("H43", "SET", 0)
    ,@("J43", "SET", 0)
    ,@("L43", "DELETE", $null)
    ,@("N43", "SET", 0)
    ,@("H46", "SET", 2145057.2)
    ,@("I46", "SET", 166.66667)
    ,@("J46", "SET", 3753725)
    ,@("K46", "SET", 500.00001)
    ,@("L46", "SET", 11261175)
    ,@("M46", "SET", -409.00001)
    ,@("N46", "SET", -11261357)
    ,@("H49", "SET", 1434.3334)
    ,@("I49", "SET", 303)
    ,@("J49", "SET", 2000)
    ,@("K49", "SET", 909)
    ,@("L49", "SET", 6000)
    ,@("M49", "SET", -753)
    ,@("N49", "SET", -6312)
    ,@("H64", "SET", 3681)
    ,@("I64", "SET", 3323.3333)
    ,@("J64", "SET", 6900)
    ,@("K64", "SET", 9969.999899999999)
    ,@("L64", "SET", 20700)
    ,@("M64", "SET", -9699.999899999999)
    ,@("N64", "SET", -21240)
    ,@("H67", "SET", 3681)
    ,@("I67", "SET", 3323.3333)
    ,@("J67", "SET", 6900)
    ,@("K67", "SET", 9969.999899999999)
    ,@("L67", "SET", 20700)
    ,@("M67", "SET", -9033.999899999999)
    ,@("N67", "SET", -22572)
    ,@("H70", "SET", 4148.294)
    ,@("I70", "SET", 2202.1)
    ,@("J70", "SET", 6928.5713)
    ,@("K70", "SET", 6606.299999999999)
    ,@("L70", "SET", 20785.7139)
    ,@("M70", "SET", -6291.299999999999)
    ,@("N70", "SET", -21415.7139)
    ,@("H73", "SET", 4148.294)
    ,@("I73", "SET", 2202.1)
    ,@("J73", "SET", 6928.5713)
    ,@("K73", "SET", 6606.299999999999)
    ,@("L73", "SET", 20785.7139)
    ,@("M73", "SET", -5514.299999999999)
    ,@("N73", "SET", -22969.7139)
    ,@("H76", "SET", 6400)
    ,@("J76", "SET", 6400)
    ,@("L76", "SET", 19200)
    ,@("N76", "SET", -19966)
    ,@("H79", "SET", 6400)
    ,@("J79", "SET", 6400)
    ,@("L79", "SET", 19200)
    ,@("N79", "SET", -21852)
    ,@("H86", "SET", 382.4)
    ,@("I86", "SET", 316.16666)
    ,@("J86", "SET", 426.55554)
    ,@("K86", "SET", 948.4999799999999)
    ,@("L86", "SET", 1279.66662)
    ,@("M86", "SET", 237.5000200000001)
    ,@("N86", "SET", -3651.66662)
    ,@("H88", "SET", 6116.6665)
    ,@("I88", "SET", 4950)
    ,@("J88", "SET", 6700)
    ,@("K88", "SET", 14850)
    ,@("L88", "SET", 20100)
    ,@("M88", "SET", -14422)
    ,@("N88", "SET", -20956)
    ,@("H89", "SET", 382.4)
    ,@("I89", "SET", 316.16666)
    ,@("J89", "SET", 426.55554)
    ,@("K89", "SET", 2845.49994)
    ,@("L89", "SET", 3838.99986)
    ,@("M89", "SET", 3082.50006)
    ,@("N89", "SET", -15694.99986)
    ,@("H91", "SET", 6116.6665)
    ,@("I91", "SET", 4950)
    ,@("J91", "SET", 6700)
    ,@("K91", "SET", 14850)
    ,@("L91", "SET", 20100)
    ,@("M91", "SET", -13368)
    ,@("N91", "SET", -23064)
    ,@("H94", "SET", 2980)
    ,@("I94", "SET", 2000)
    ,@("J94", "SET", 6900)
    ,@("K94", "SET", 6000)
    ,@("L94", "SET", 20700)
    ,@("M94", "SET", -5324)
    ,@("N94", "SET", -22052)
    ,@("H123", "SET", 4756.143)
    ,@("J123", "SET", 7433.25)
    ,@("L123", "SET", 22299.75)
    ,@("N123", "SET", -27199.75)
)
$changes["GSM"] = @(
    ,@("H2", "SET", 48.333332)
    ,@("I2", "SET", 27.333334)
    ,@("J2", "SET", 90.333336)
    ,@("K2", "SET", 27.333334)
    ,@("L2", "SET", 90.333336)
    ,@("M2", "SET", 85.66666599999999)
    ,@("N2", "SET", -316.333336)
)
$changes["LTW"] = @(
    ,@("H61", "SET", 1836.7142)
    ,@("I61", "SET", 1371.9)
    ,@("J61", "SET", 2998.75)
    ,@("K61", "SET", 1371.9)
    ,@("L61", "SET", 2998.75)
    ,@("M61", "SET", -1169.9)
    ,@("N61", "SET", -3402.75)
    ,@("H113", "SET", 1836.7142)
    ,@("I113", "SET", 1371.9)
    ,@("J113", "SET", 2998.75)
    ,@("K113", "SET", 1371.9)
    ,@("L113", "SET", 2998.75)
    ,@("M113", "SET", 798.0999999999999)
    ,@("N113", "SET", -7338.75)
    ,@("H132", "SET", 3513.6667)
    ,@("I132", "SET", 2088.375)
    ,@("J132", "SET", 5142.5713)
    ,@("K132", "SET", 6265.125)
    ,@("L132", "SET", 15427.7139)
    ,@("M132", "SET", -3735.125)
    ,@("N132", "SET", -20487.7139)
)
$changes["WVR"] = @(
    ,@("H86", "SET", 29966.666)
    ,@("J86", "SET", 29966.666)
    ,@("L86", "SET", 29966.666)
    ,@("N86", "SET", -32212.666)
    ,@("H89", "SET", 29966.666)
    ,@("J89", "SET", 29966.666)
    ,@("L89", "SET", 149833.33)
    ,@("N89", "SET", -161065.33)
    ,@("H113", "SET", 204.8)
    ,@("I113", "SET", 193.88889)
    ,@("K113", "SET", 581.6666700000001)
    ,@("M113", "SET", 1588.33333)
)

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($change in $changes[$sheetName]) {
        $cellRef = $change[0]
        $action = $change[1]
        $value = $change[2]
        if ($action -eq "SET") {
            $ws.Range($cellRef).Value = $value
        } else {
            $ws.Range($cellRef).ClearContents()
        }
    }
}

Write-Output "Applied $(($changes.Values | ForEach-Object { $_.Count } | Measure-Object -Sum).Sum) cell changes across $($changes.Keys.Count) sheets"